$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new numeric index column (A2:A21) the same formatting/style
# as the existing header cells (bold, bordered, centered) by copying
# the format from B1 and pasting only formats into A2:A21.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A2:A21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 'RPA Developer'
$ws.Cells.Item(2, 3).Value = 'https://www.naukri.com/job-listings-rpa-developer-dautom-bengaluru-2-to-6-years-170524500436'

$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 'Automation Anywhere Developer'
$ws.Cells.Item(3, 3).Value = 'https://www.naukri.com/job-listings-automation-anywhere-developer-ilink-digital-kolkata-mumbai-new-delhi-hyderabad-pune-chennai-bengaluru-2-to-4-years-170524500687'

$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 'RPA Business Analyst'
$ws.Cells.Item(4, 3).Value = 'https://www.naukri.com/job-listings-rpa-business-analyst-percipere-mumbai-pune-0-to-1-years-170524909865'

$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = 'Senior RPA Business Analyst'
$ws.Cells.Item(5, 3).Value = 'https://www.naukri.com/job-listings-senior-rpa-business-analyst-percipere-mumbai-5-to-7-years-170524909645'

$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = 'Rpa Developer'
$ws.Cells.Item(6, 3).Value = 'https://www.naukri.com/job-listings-rpa-developer-apmosys-technologies-navi-mumbai-8-to-10-years-170524008223'

$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = 'Opportunity For RPA Developers'
$ws.Cells.Item(7, 3).Value = 'https://www.naukri.com/job-listings-opportunity-for-rpa-developers-healthcare-informatics-vadodara-2-to-5-years-250424004964'

$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = 'RPA Support Specialist (Night Shift, US Hours)'
$ws.Cells.Item(8, 3).Value = 'https://www.naukri.com/job-listings-rpa-support-specialist-night-shift-us-hours-percipere-mumbai-2-to-3-years-170524912102'

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = 'RPA diversity opening For Ahmedabad Location'
$ws.Cells.Item(9, 3).Value = 'https://www.naukri.com/job-listings-rpa-diversity-opening-for-ahmedabad-location-infosys-bpm-ahmedabad-6-to-10-years-030424012572'

$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = 'Senior Developer'
$ws.Cells.Item(10, 3).Value = 'https://www.naukri.com/job-listings-senior-developer-acronotics-pvt-ltd-pune-bengaluru-5-to-10-years-170524500545'

$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = 'Developer'
$ws.Cells.Item(11, 3).Value = 'https://www.naukri.com/job-listings-developer-acronotics-pvt-ltd-pune-bengaluru-3-to-8-years-170524500544'

$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = 'Solution Architects'
$ws.Cells.Item(12, 3).Value = 'https://www.naukri.com/job-listings-solution-architects-acronotics-pvt-ltd-pune-bengaluru-8-to-13-years-170524500638'

$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = 'RPA  Robot Framework & Python QA Automation'
$ws.Cells.Item(13, 3).Value = 'https://www.naukri.com/job-listings-rpa-robot-framework-python-qa-automation-rq-technologies-llp-chennai-5-to-10-years-170524011474'

$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = 'Application Developer'
$ws.Cells.Item(14, 3).Value = 'https://www.naukri.com/job-listings-application-developer-accenture-solutions-pvt-ltd-mumbai-3-to-5-years-170524908871'

$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = 'Application Developer'
$ws.Cells.Item(15, 3).Value = 'https://www.naukri.com/job-listings-application-developer-accenture-solutions-pvt-ltd-chennai-3-to-6-years-170524904578'

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = 'Technical BotOps'
$ws.Cells.Item(16, 3).Value = 'https://www.naukri.com/job-listings-technical-botops-allegis-group-hyderabad-1-to-5-years-170524004542'

$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = 'Application Designer'
$ws.Cells.Item(17, 3).Value = 'https://www.naukri.com/job-listings-application-designer-accenture-solutions-pvt-ltd-bengaluru-3-to-5-years-170524912002'

$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = 'Application Designer'
$ws.Cells.Item(18, 3).Value = 'https://www.naukri.com/job-listings-application-designer-accenture-solutions-pvt-ltd-bengaluru-3-to-7-years-170524911413'

$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = 'Application Designer'
$ws.Cells.Item(19, 3).Value = 'https://www.naukri.com/job-listings-application-designer-accenture-solutions-pvt-ltd-bengaluru-7-to-9-years-170524912884'

$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 2).Value = 'Trust & Safety New Associate'
$ws.Cells.Item(20, 3).Value = 'https://www.naukri.com/job-listings-trust-safety-new-associate-accenture-solutions-pvt-ltd-gurugram-0-to-1-years-180524908538'

$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).Value = 'S&C Global Network - Strategy - MC - Industry X - Capital Projects'
$ws.Cells.Item(21, 3).Value = 'https://www.naukri.com/job-listings-s-c-global-network-strategy-mc-industry-x-capital-projects-accenture-solutions-pvt-ltd-gurugram-7-to-9-years-170524908957'
